$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("CAB:Station:958219cd-8009-45b6-8a04-1cc4ae763307", 5, 1709651610, "2024-03-05 16:16:16.319954"),
    @("CAB:Station:69a20737-5c07-4ab4-ad8a-4fcba602ef6b", 3, 1709651610, "2024-03-05 16:16:16.319954"),
    @("CAB:Station:59830a25-c7c6-4c9b-af21-0b67c836ba8f", 5, 1709651610, "2024-03-05 16:16:16.319954"),
    @("CAB:Station:36fe6a30-d53e-4aa8-863d-e21e80fd1d0b", 5, 1709651610, "2024-03-05 16:16:16.319954"),
    @("CAB:Station:ca41b63e-9046-442f-880e-d0e9186e507a", 1, 1709651610, "2024-03-05 16:16:16.319954"),
    @("CAB:Station:714932fd-4836-4a2a-9a19-07a6d694274c", 7, 1709651610, "2024-03-05 16:16:16.319954"),
    @("CAB:Station:7046d131-47fd-40bd-ba63-f52590d2ee8f", 2, 1709651610, "2024-03-05 16:16:16.319954"),
    @("CAB:Station:e375c577-ca9b-4f3d-b052-582f4eeeb0a8", 2, 1709651610, "2024-03-05 16:16:16.319954"),
    @("CAB:Station:f384359e-349d-47cf-aec3-031f20b8e0a0", 2, 1709651610, "2024-03-05 16:16:16.319954"),
    @("CAB:Station:e73ff873-8f43-41b4-868a-2d3e00f2cfbb", 2, 1709651610, "2024-03-05 16:16:16.319954"),
    @("CAB:Station:c070e40c-98f6-45d6-9543-90ceceef63af", 9, 1709651610, "2024-03-05 16:16:16.319954"),
    @("CAB:Station:2263b695-db35-4bbf-b70a-032920fadf3f", 1, 1709651610, "2024-03-05 16:16:16.319954"),
    @("CAB:Station:be796eb3-49af-494c-b3f0-dc1f0132eb7d", 2, 1709651610, "2024-03-05 16:16:16.319954"),
    @("CAB:Station:b207ab7a-f23d-40f8-aa26-5d02077fc6f9", 2, 1709651610, "2024-03-05 16:16:16.319954")
)

$startRow = 59
for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $startRow + $i
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
